$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old contents for columns B:H in row 1 and rows 2:3 entirely (B:H),
# and remove row 3 altogether, leaving only A1:A2 used.
$ws.Range("B1:H3").Clear()
$ws.Rows.Item(3).Delete()

# Update remaining cells
$ws.Range("A1").Value = "intervalo"
$ws.Range("A2").Value = -1
